$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B2').Value = 'hold1'
$ws.Range('C2').Value = 'møk2'
$ws.Range('D2').Value = 'TK'
$ws.Range('E2').Value = 'fys5'
$ws.Range('B3').Value = 'fys4'
$ws.Range('C3').Value = 'hold2'
$ws.Range('D3').Value = 'dat5'
$ws.Range('E3').Value = 'fys2'
$ws.Range('B4').Value = 'mat2'
$ws.Range('C4').Value = 'møk1'
$ws.Range('D4').Value = 'dat3'
$ws.Range('E4').Value = 'dat1'
$ws.Range('B5').Value = 'dat6'
$ws.Range('C5').Value = 'dat7'
$ws.Range('D5').Value = 'dat4'
$ws.Range('E5').Value = 'fys1'
$ws.Range('B6').Value = 'nano'
$ws.Range('C6').Value = 'dav'
$ws.Range('D6').Value = 'mat3'
$ws.Range('E6').Value = 'it1'
$ws.Range('B7').Value = 'mat1'
$ws.Range('C7').Value = 'dat2'
$ws.Range('D7').Value = 'it2'
$ws.Range('E7').Value = 'fys3'
$ws.Range('B11').Value = 'hold2'
$ws.Range('C11').Value = 'dat3'
$ws.Range('D11').Value = 'dav'
$ws.Range('B12').Value = 'hold1'
$ws.Range('D12').Value = 'dat4'
$ws.Range('E12').Value = 'mat1'
$ws.Range('B13').Value = 'fys5'
$ws.Range('C13').Value = 'fys4'
$ws.Range('E13').Value = 'fys3'
$ws.Range('B14').Value = 'TK'
$ws.Range('C14').Value = 'dat5'
$ws.Range('D14').Value = 'mat2'
$ws.Range('E14').Value = 'mat3'
$ws.Range('B15').Value = 'møk2'
$ws.Range('C15').Value = 'fys2'
$ws.Range('D15').Value = 'dat1'
$ws.Range('B16').Value = 'dat6'
$ws.Range('C16').Value = 'dat7'
$ws.Range('D16').Value = 'nano'
$ws.Range('E16').Value = 'it1'
$ws.Range('B20').Value = 'dat5'
$ws.Range('C20').Value = 'dat1'
$ws.Range('D20').Value = 'dat6'
$ws.Range('E20').Value = 'it2'
$ws.Range('B21').Value = 'fys5'
$ws.Range('C21').Value = 'fys1'
$ws.Range('D21').Value = 'nano'
$ws.Range('B22').Value = 'hold1'
$ws.Range('C22').Value = 'fys2'
$ws.Range('D22').Value = 'dat7'
$ws.Range('E22').Value = 'it1'
$ws.Range('B23').Value = 'møk2'
$ws.Range('C23').Value = 'hold2'
$ws.Range('D23').Value = 'mat1'
$ws.Range('E23').Value = 'fys3'
$ws.Range('B24').Value = 'TK'
$ws.Range('C24').Value = 'fys4'
$ws.Range('D24').Value = 'dat3'
$ws.Range('E24').Value = 'dat4'
$ws.Range('B25').Value = 'mat2'
$ws.Range('C25').Value = 'møk1'
$ws.Range('D25').Value = 'dav'
$ws.Range('E25').Value = 'mat3'
$ws.Range('B29').Value = 'møk1'
$ws.Range('C29').Value = 'fys1'
$ws.Range('D29').Value = 'it1'
$ws.Range('E29').Value = 'fys3'
$ws.Range('B30').Value = 'TK'
$ws.Range('C30').Value = 'dat1'
$ws.Range('D30').Value = 'dat7'
$ws.Range('E30').Value = 'dav'
$ws.Range('B31').Value = 'møk2'
$ws.Range('C31').Value = 'dat4'
$ws.Range('D31').Value = 'mat3'
$ws.Range('B32').Value = 'hold1'
$ws.Range('C32').Value = 'dat3'
$ws.Range('D32').Value = 'nano'
$ws.Range('E32').Value = 'it2'
$ws.Range('B33').Value = 'fys5'
$ws.Range('D33').Value = 'dat6'
$ws.Range('B34').Value = 'fys4'
$ws.Range('C34').Value = 'hold2'
$ws.Range('D34').Value = 'dat5'
$ws.Range('E34').Value = 'fys2'
$ws.Range('B38').Value = 'fys4'
$ws.Range('C38').Value = 'dat7'
$ws.Range('D38').Value = 'mat3'
$ws.Range('E38').Value = 'mat1'
$ws.Range('B39').Value = 'møk2'
$ws.Range('C39').Value = 'mat2'
$ws.Range('D39').Value = 'it1'
$ws.Range('E39').Value = 'it2'
$ws.Range('B40').Value = 'TK'
$ws.Range('C40').Value = 'hold2'
$ws.Range('D40').Value = 'dat6'
$ws.Range('E40').Value = 'nano'
$ws.Range('B41').Value = 'fys5'
$ws.Range('C41').Value = 'fys2'
$ws.Range('D41').Value = 'møk1'
$ws.Range('E41').Value = 'dav'
$ws.Range('B42').Value = 'hold1'
$ws.Range('C42').Value = 'dat5'
$ws.Range('D42').Value = 'dat2'
$ws.Range('E42').Value = 'fys3'
$ws.Range('B43').Value = 'dat3'
$ws.Range('C43').Value = 'dat1'
$ws.Range('D43').Value = 'dat4'
$ws.Range('E43').Value = 'fys1'
$ws.Range('B47').Value = 'fys2'
$ws.Range('C47').Value = 'mat2'
$ws.Range('D47').Value = 'dat4'
$ws.Range('E47').Value = 'nano'
$ws.Range('B48').Value = 'dat3'
$ws.Range('C48').Value = 'dat6'
$ws.Range('D48').Value = 'mat3'
$ws.Range('E48').Value = 'fys3'
$ws.Range('B49').Value = 'dat5'
$ws.Range('C49').Value = 'fys1'
$ws.Range('D49').Value = 'dav'
$ws.Range('E49').Value = 'mat1'
$ws.Range('B50').Value = 'fys4'
$ws.Range('C50').Value = 'dat1'
$ws.Range('D50').Value = 'it1'
$ws.Range('B51').Value = 'hold2'
$ws.Range('C51').Value = 'møk1'
$ws.Range('D51').Value = 'dat7'
$ws.Range('E51').Value = 'it2'
$ws.Range('B52').Value = 'hold1'
$ws.Range('C52').Value = 'møk2'
$ws.Range('D52').Value = 'TK'
$ws.Range('E52').Value = 'fys5'
